# Applies the ARP_CAM_System.pptx diff on slide 2 (the thread/process diagram):
#   - Resizes the outer background parallelogram.
#   - Resizes/repositions the four inner parallelograms (threads) and their
#     caption textboxes.
#   - Renames the generic "Thread" captions to the specific thread names.
#
# Note: PowerPoint's Shape.Left/Top/Width/Height are single-precision (Single)
# floats measured in points (1 pt = 12700 EMU). The point literals below were
# chosen so that, after the host's float32 round-trip, they resolve back to
# the exact target EMU values from the source XML (plain emu/12700 can land
# 1 EMU short because of float32 truncation).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Shape 1: "Parallelogram 2" - outer/background parallelogram
$shp = $s.Shapes.Item(1)
$shp.Left   = 94.00023622047244
$shp.Top    = 100.96299412598425
$shp.Width  = 450.4441732283465
$shp.Height = 264.23692913385827

# Shape 3: "Parallelogram 4"
$shp = $s.Shapes.Item(3)
$shp.Left   = 146.7164566929134
$shp.Top    = 247.62118110236221
$shp.Width  = 144.0
$shp.Height = 72.0

# Shape 4: "TextBox 9" -> "thread get_attitude"
$shp = $s.Shapes.Item(4)
$shp.Left   = 154.63645669291338
$shp.Top    = 265.6211811023622
$shp.Width  = 128.16
$shp.Height = 36.0
$shp.TextFrame.TextRange.Text = "thread get_attitude"

# Shape 5: "Parallelogram 10"
$shp = $s.Shapes.Item(5)
$shp.Left   = 170.28740757480315
$shp.Top    = 151.8
$shp.Width  = 144.0
$shp.Height = 72.0

# Shape 6: "TextBox 11" -> "thread get_radio_command"
$shp = $s.Shapes.Item(6)
$shp.Left   = 178.20740557480315
$shp.Top    = 169.8
$shp.Width  = 128.16
$shp.Height = 36.0
$shp.TextFrame.TextRange.Text = "thread get_radio_command"

# Shape 7: "Parallelogram 12"
$shp = $s.Shapes.Item(7)
$shp.Left   = 343.81795275590554
$shp.Top    = 151.8
$shp.Width  = 144.0
$shp.Height = 72.0

# Shape 8: "TextBox 13" -> "thread control_servos"
$shp = $s.Shapes.Item(8)
$shp.Left   = 351.7379617559055
$shp.Top    = 169.8
$shp.Width  = 128.16
$shp.Height = 36.0
$shp.TextFrame.TextRange.Text = "thread control_servos"

# Shape 9: "Parallelogram 14"
$shp = $s.Shapes.Item(9)
$shp.Left   = 319.30946381889765
$shp.Top    = 247.62118110236221
$shp.Width  = 144.0
$shp.Height = 72.0

# Shape 10: "TextBox 15" -> "thread process_video"
$shp = $s.Shapes.Item(10)
$shp.Left   = 327.22944881889765
$shp.Top    = 265.6211811023622
$shp.Width  = 128.16
$shp.Height = 36.0
$shp.TextFrame.TextRange.Text = "thread process_video"
